$wb = $excel.ActiveWorkbook

# --- Sheet "Главные" (Main officials) ---
$wsMain = $wb.Worksheets.Item("Главные")

# Row 2
$wsMain.Range("C2").Value = 23
$wsMain.Range("D2").Value = 530
$wsMain.Range("E2").Value = 230
$wsMain.Range("F2").Value = 300
$wsMain.Range("G2").Value = 23.04
$wsMain.Range("H2").Value = 10
$wsMain.Range("I2").Value = 13.04
$wsMain.Range("J2").Value = 100
$wsMain.Range("K2").Value = 120

# Row 3
$wsMain.Range("C3").Value = 23
$wsMain.Range("D3").Value = 417
$wsMain.Range("E3").Value = 183
$wsMain.Range("F3").Value = 234
$wsMain.Range("G3").Value = 18.13
$wsMain.Range("H3").Value = 7.96
$wsMain.Range("I3").Value = 10.17
$wsMain.Range("J3").Value = 89
$wsMain.Range("K3").Value = 97
$wsMain.Range("V3").Value = 12

# Row 6
$wsMain.Range("C6").Value = 23
$wsMain.Range("D6").Value = 411
$wsMain.Range("E6").Value = 171
$wsMain.Range("F6").Value = 240
$wsMain.Range("G6").Value = 17.87
$wsMain.Range("H6").Value = 7.43
$wsMain.Range("I6").Value = 10.43
$wsMain.Range("J6").Value = 78
$wsMain.Range("K6").Value = 100

# Row 7
$wsMain.Range("C7").Value = 15
$wsMain.Range("D7").Value = 203
$wsMain.Range("E7").Value = 88
$wsMain.Range("F7").Value = 115
$wsMain.Range("G7").Value = 13.53
$wsMain.Range("H7").Value = 5.87
$wsMain.Range("I7").Value = 7.67
$wsMain.Range("J7").Value = 44
$wsMain.Range("K7").Value = 40
$wsMain.Range("V7").Value = 12

# Row 24
$wsMain.Range("C24").Value = 23
$wsMain.Range("D24").Value = 364
$wsMain.Range("E24").Value = 173
$wsMain.Range("F24").Value = 191
$wsMain.Range("G24").Value = 15.83
$wsMain.Range("H24").Value = 7.52
$wsMain.Range("I24").Value = 8.3
$wsMain.Range("J24").Value = 84
$wsMain.Range("K24").Value = 93

# Row 25
$wsMain.Range("C25").Value = 23
$wsMain.Range("D25").Value = 390
$wsMain.Range("E25").Value = 196
$wsMain.Range("F25").Value = 194
$wsMain.Range("G25").Value = 16.96
$wsMain.Range("H25").Value = 8.52
$wsMain.Range("I25").Value = 8.43
$wsMain.Range("J25").Value = 93
$wsMain.Range("K25").Value = 92

# --- Sheet "Линейные" (Linesmen) ---
$wsLine = $wb.Worksheets.Item("Линейные")

# Row 8
$wsLine.Range("C8").Value = 21
$wsLine.Range("D8").Value = 312
$wsLine.Range("E8").Value = 128
$wsLine.Range("F8").Value = 184
$wsLine.Range("G8").Value = 14.86
$wsLine.Range("H8").Value = 6.1
$wsLine.Range("I8").Value = 8.76
$wsLine.Range("J8").Value = 59
$wsLine.Range("K8").Value = 77

# Row 9
$wsLine.Range("C9").Value = 21
$wsLine.Range("D9").Value = 395
$wsLine.Range("E9").Value = 181
$wsLine.Range("F9").Value = 214
$wsLine.Range("G9").Value = 18.81
$wsLine.Range("H9").Value = 8.62
$wsLine.Range("I9").Value = 10.19
$wsLine.Range("J9").Value = 78
$wsLine.Range("K9").Value = 97

# Row 11
$wsLine.Range("C11").Value = 14
$wsLine.Range("D11").Value = 207
$wsLine.Range("E11").Value = 100
$wsLine.Range("F11").Value = 107
$wsLine.Range("G11").Value = 14.79
$wsLine.Range("H11").Value = 7.14
$wsLine.Range("I11").Value = 7.64
$wsLine.Range("J11").Value = 50
$wsLine.Range("K11").Value = 51
$wsLine.Range("V11").Value = 6

# Row 16
$wsLine.Range("C16").Value = 22
$wsLine.Range("D16").Value = 376
$wsLine.Range("E16").Value = 182
$wsLine.Range("F16").Value = 194
$wsLine.Range("G16").Value = 17.09
$wsLine.Range("H16").Value = 8.27
$wsLine.Range("I16").Value = 8.82
$wsLine.Range("J16").Value = 86
$wsLine.Range("K16").Value = 92
$wsLine.Range("V16").Value = 12

# Row 17
$wsLine.Range("C17").Value = 11
$wsLine.Range("D17").Value = 150
$wsLine.Range("E17").Value = 88
$wsLine.Range("F17").Value = 62
$wsLine.Range("G17").Value = 13.64
$wsLine.Range("H17").Value = 8
$wsLine.Range("I17").Value = 5.64
$wsLine.Range("J17").Value = 44
$wsLine.Range("K17").Value = 31

# Row 18
$wsLine.Range("C18").Value = 24
$wsLine.Range("D18").Value = 409
$wsLine.Range("E18").Value = 194
$wsLine.Range("F18").Value = 215
$wsLine.Range("G18").Value = 17.04
$wsLine.Range("H18").Value = 8.08
$wsLine.Range("I18").Value = 8.96
$wsLine.Range("J18").Value = 92
$wsLine.Range("K18").Value = 90

# --- Update "as_of_utc" timestamp (column AA) for all data rows (2-26) on both sheets ---
$newTimestamp = "2025-11-10 03:06:39"
for ($r = 2; $r -le 26; $r++) {
    $wsMain.Range("AA$r").Value = $newTimestamp
    $wsLine.Range("AA$r").Value = $newTimestamp
}

Write-Host "Applied KHL referees stats update (2025-11-10 03:06:39)"
